# Update "想去人数" (interest count) values in column F across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.
# 本地生活 (Local life) sheet is unchanged.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1014
$ws1.Range("F6").Value  = 5678
$ws1.Range("F8").Value  = 743
$ws1.Range("F9").Value  = 979
$ws1.Range("F19").Value = 986
$ws1.Range("F23").Value = 580
$ws1.Range("F24").Value = 176
$ws1.Range("F28").Value = 3171
$ws1.Range("F31").Value = 76
$ws1.Range("F32").Value = 138
$ws1.Range("F34").Value = 433
$ws1.Range("F39").Value = 306
$ws1.Range("F40").Value = 764
$ws1.Range("F44").Value = 74

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 224
$ws2.Range("F6").Value = 156

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1014
$ws4.Range("F6").Value  = 5678
$ws4.Range("F8").Value  = 743
$ws4.Range("F10").Value = 224
$ws4.Range("F11").Value = 979
$ws4.Range("F14").Value = 156
$ws4.Range("F24").Value = 986
$ws4.Range("F28").Value = 580
$ws4.Range("F29").Value = 176
$ws4.Range("F31").Value = 3171
$ws4.Range("F34").Value = 76
$ws4.Range("F35").Value = 138
$ws4.Range("F37").Value = 433
$ws4.Range("F41").Value = 306
$ws4.Range("F42").Value = 764
$ws4.Range("F46").Value = 74
